# Scheduled-runner style update of cached market-price figures across the
# per-job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). These are
# plain cached values (no formulas), so each touched cell is written
# directly via Range.Value; ARM!N122 loses its HQ-profit figure entirely
# once the HQ price drops to 0, so it is cleared instead of being set.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 616.625
$ws.Range("I28").Value = 597.7143
$ws.Range("K28").Value = 597.7143
$ws.Range("M28").Value = -112.7143
$ws.Range("H33").Value = 6015.857
$ws.Range("I33").Value = 7857.769
$ws.Range("K33").Value = 7857.769
$ws.Range("M33").Value = -7628.769
$ws.Range("H43").Value = 55573640
$ws.Range("J43").Value = 26623.75
$ws.Range("L43").Value = 26623.75
$ws.Range("N43").Value = -26761.75
$ws.Range("H88").Value = 5369
$ws.Range("I88").Value = 1291.6666
$ws.Range("J88").Value = 8427
$ws.Range("K88").Value = 1291.6666
$ws.Range("L88").Value = 8427
$ws.Range("M88").Value = -885.6666
$ws.Range("N88").Value = -9239
$ws.Range("H91").Value = 5369
$ws.Range("I91").Value = 1291.6666
$ws.Range("J91").Value = 8427
$ws.Range("K91").Value = 1291.6666
$ws.Range("L91").Value = 8427
$ws.Range("M91").Value = 112.3334
$ws.Range("N91").Value = -11235
$ws.Range("H123").Value = 136330
$ws.Range("J123").Value = 136330
$ws.Range("L123").Value = 136330
$ws.Range("N123").Value = -146130
$ws.Range("H141").Value = 5549.4443
$ws.Range("I141").Value = 5618.875
$ws.Range("K141").Value = 16856.625
$ws.Range("M141").Value = -11676.625

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 5679
$ws.Range("I46").Value = 1999.5
$ws.Range("K46").Value = 1999.5
$ws.Range("M46").Value = -1680.5
$ws.Range("H97").Value = 5303.6665
$ws.Range("I97").Value = 2950
$ws.Range("K97").Value = 2950
$ws.Range("M97").Value = -2454
$ws.Range("H102").Value = 3255
$ws.Range("I102").Value = 3255
$ws.Range("K102").Value = 3255
$ws.Range("M102").Value = -1633
$ws.Range("H122").Value = 3188.2942
$ws.Range("I122").Value = 3188.2942
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9564.882599999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7114.882599999999
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 3185.4167
$ws.Range("I132").Value = 2518.889
$ws.Range("K132").Value = 7556.667
$ws.Range("M132").Value = -5026.667

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 62503344
$ws.Range("J80").Value = 3793.2727
$ws.Range("L80").Value = 3793.2727
$ws.Range("N80").Value = -5789.2727
$ws.Range("H83").Value = 62503344
$ws.Range("J83").Value = 3793.2727
$ws.Range("L83").Value = 18966.3635
$ws.Range("N83").Value = -28950.3635
$ws.Range("H86").Value = 4962.9287
$ws.Range("I86").Value = 2221.889
$ws.Range("K86").Value = 2221.889
$ws.Range("M86").Value = -1098.889
$ws.Range("H89").Value = 4962.9287
$ws.Range("I89").Value = 2221.889
$ws.Range("K89").Value = 11109.445
$ws.Range("M89").Value = -5493.445
$ws.Range("H94").Value = 2827.7144
$ws.Range("I94").Value = 1611.9333
$ws.Range("K94").Value = 1611.9333
$ws.Range("M94").Value = -1160.9333
$ws.Range("H99").Value = 12901
$ws.Range("I99").Value = 13990.7
$ws.Range("K99").Value = 13990.7
$ws.Range("M99").Value = -12492.7
$ws.Range("H105").Value = 22782.75
$ws.Range("I105").Value = 7899.25
$ws.Range("K105").Value = 7899.25
$ws.Range("M105").Value = -6152.25
$ws.Range("H134").Value = 37502636
$ws.Range("I134").Value = 2331.25
$ws.Range("J134").Value = 112503250
$ws.Range("K134").Value = 6993.75
$ws.Range("L134").Value = 337509750
$ws.Range("M134").Value = -4458.75
$ws.Range("N134").Value = -337514820

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2201.8096
$ws.Range("I58").Value = 2420.8
$ws.Range("K58").Value = 2420.8
$ws.Range("M58").Value = -2217.8
$ws.Range("H132").Value = 3386
$ws.Range("I132").Value = 3441.625
$ws.Range("K132").Value = 10324.875
$ws.Range("M132").Value = -7794.875
$ws.Range("H136").Value = 2201.8096
$ws.Range("I136").Value = 2420.8
$ws.Range("K136").Value = 7262.400000000001
$ws.Range("M136").Value = -4712.400000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1560.5
$ws.Range("I5").Value = 1156.8889
$ws.Range("J5").Value = 1890.7273
$ws.Range("K5").Value = 3470.6667
$ws.Range("L5").Value = 5672.1819
$ws.Range("M5").Value = -3358.6667
$ws.Range("N5").Value = -5896.1819
$ws.Range("H12").Value = 275.57144
$ws.Range("J12").Value = 328
$ws.Range("L12").Value = 984
$ws.Range("N12").Value = -1330
$ws.Range("H135").Value = 1560.5
$ws.Range("I135").Value = 1156.8889
$ws.Range("J135").Value = 1890.7273
$ws.Range("K135").Value = 10412.0001
$ws.Range("L135").Value = 17016.5457
$ws.Range("M135").Value = -7877.000099999999
$ws.Range("N135").Value = -22086.5457

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 41668216
$ws.Range("I102").Value = 41668216
$ws.Range("K102").Value = 41668216
$ws.Range("M102").Value = -41666594
$ws.Range("H132").Value = 856257.0600000001
$ws.Range("I132").Value = 927.7143
$ws.Range("J132").Value = 1116574.8
$ws.Range("K132").Value = 2783.1429
$ws.Range("L132").Value = 3349724.4
$ws.Range("M132").Value = -253.1428999999998
$ws.Range("N132").Value = -3354784.4
$ws.Range("H137").Value = 183138.67
$ws.Range("J137").Value = 189759
$ws.Range("L137").Value = 189759
$ws.Range("N137").Value = -199959

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10444.75
$ws.Range("J7").Value = 19398.6
$ws.Range("L7").Value = 19398.6
$ws.Range("N7").Value = -19622.6
$ws.Range("H55").Value = 639.5135
$ws.Range("J55").Value = 809.63635
$ws.Range("L55").Value = 809.63635
$ws.Range("N55").Value = -1155.63635
$ws.Range("H82").Value = 1646.0714
$ws.Range("I82").Value = 1616.091
$ws.Range("J82").Value = 1756
$ws.Range("K82").Value = 1616.091
$ws.Range("L82").Value = 1756
$ws.Range("M82").Value = -1255.091
$ws.Range("N82").Value = -2478
$ws.Range("H85").Value = 1646.0714
$ws.Range("I85").Value = 1616.091
$ws.Range("J85").Value = 1756
$ws.Range("K85").Value = 1616.091
$ws.Range("L85").Value = 1756
$ws.Range("M85").Value = -368.0909999999999
$ws.Range("N85").Value = -4252
$ws.Range("H100").Value = 4724.7334
$ws.Range("I100").Value = 4187.25
$ws.Range("J100").Value = 4920.1816
$ws.Range("K100").Value = 4187.25
$ws.Range("L100").Value = 4920.1816
$ws.Range("M100").Value = -3646.25
$ws.Range("N100").Value = -6002.1816
$ws.Range("H126").Value = 10444.75
$ws.Range("J126").Value = 19398.6
$ws.Range("L126").Value = 58195.8
$ws.Range("N126").Value = -63135.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 860.5
$ws.Range("I100").Value = 661.0714
$ws.Range("K100").Value = 1322.1428
$ws.Range("M100").Value = -781.1428000000001
$ws.Range("H119").Value = 54229.668
$ws.Range("J119").Value = 54229.668
$ws.Range("L119").Value = 54229.668
$ws.Range("N119").Value = -63905.668
$ws.Range("H123").Value = 84993.39999999999
$ws.Range("J123").Value = 84993.39999999999
$ws.Range("L123").Value = 84993.39999999999
$ws.Range("N123").Value = -94793.39999999999
$ws.Range("H126").Value = 2682.2727
$ws.Range("I126").Value = 2586.842
$ws.Range("K126").Value = 7760.526
$ws.Range("M126").Value = -5290.526
$ws.Range("H140").Value = 79999.25
$ws.Range("J140").Value = 79999.25
$ws.Range("L140").Value = 79999.25
$ws.Range("N140").Value = -90359.25
